# Updated cryptos list on Sat Jan  6 14:25:46 UTC 2024 with GitHub Actions
#
# Refreshes the scraped coinranking.com snapshot on Sheet1 (columns: A=rank,
# B=Coin, C=Link, D=Price, E=Volume(1h)). Re-applies each changed cell from the
# latest scrape; a handful of rows also re-rank (name/link/price/volume all move
# together) where coins swapped places in the source ranking.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '44.121.45'
$ws.Range("E2").Value = '  +0.44%  '

# Row 3
$ws.Range("D3").Value = '2.245.62'
$ws.Range("E3").Value = '  +0.37%  '

# Row 4
$ws.Range("E4").Value = '  -0.02%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '306.70'
$ws.Range("E5").Value = '  -2.81%  '
$ws.Range("D5").Style = "Normal"

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '96.24'
$ws.Range("E6").Value = '  -3.36%  '
$ws.Range("D6").Style = "Normal"

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.573'
$ws.Range("E7").Value = '  +0.46%  '
$ws.Range("D7").Style = "Normal"

# Row 8
$ws.Range("E8").Value = '  +0.15%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.528'
$ws.Range("E9").Value = '  -1.62%  '
$ws.Range("D9").Style = "Normal"

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.07'
$ws.Range("E10").Value = '  -3.73%  '
$ws.Range("D10").Style = "Normal"

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0812'
$ws.Range("E11").Value = '  -1.03%  '
$ws.Range("D11").Style = "Normal"

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.25'
$ws.Range("E12").Value = '  -1.75%  '
$ws.Range("D12").Style = "Normal"

# Row 13
$ws.Range("E13").Value = '  -0.35%  '

# Row 14
$ws.Range("D14").Value = '2.587.06'
$ws.Range("E14").Value = '  +0.17%  '

# Row 15
$ws.Range("D15").Value = '2.253.08'
$ws.Range("E15").Value = '  +0.55%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.829'
$ws.Range("E16").Value = '  -1.51%  '
$ws.Range("D16").Style = "Normal"

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '13.61'
$ws.Range("E17").Value = '  -2.99%  '
$ws.Range("D17").Style = "Normal"

# Row 18
$ws.Range("D18").Value = '44.039.19'
$ws.Range("E18").Value = '  +0.39%  '

# Row 19
$ws.Range("D19").Value = '0.0₃0969'
$ws.Range("E19").Value = '  +0.40%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.30'
$ws.Range("E20").Value = '  -4.87%  '
$ws.Range("D20").Style = "Normal"

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.33'
$ws.Range("E21").Value = '  -0.40%  '
$ws.Range("D21").Style = "Normal"

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '65.57'
$ws.Range("E22").Value = '  +0.89%  '
$ws.Range("D22").Style = "Normal"

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '236.73'
$ws.Range("E23").Value = '  +1.17%  '
$ws.Range("D23").Style = "Normal"

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.94'
$ws.Range("E24").Value = '  -3.28%  '
$ws.Range("D24").Style = "Normal"

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.99'
$ws.Range("E25").Value = '  -2.80%  '
$ws.Range("D25").Style = "Normal"

# Row 26
$ws.Range("E26").Value = '  +0.00%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '38.85'
$ws.Range("E27").Value = '  +5.61%  '
$ws.Range("D27").Style = "Normal"

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.92'
$ws.Range("E28").Value = '  -3.42%  '
$ws.Range("D28").Style = "Normal"

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.21'
$ws.Range("E29").Value = '  +0.99%  '
$ws.Range("D29").Style = "Normal"

# Row 30
$ws.Range("E30").Value = '  +0.79%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '20.03'
$ws.Range("E31").Value = '  +0.40%  '
$ws.Range("D31").Style = "Normal"

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '151.63'
$ws.Range("E32").Value = '  -4.46%  '
$ws.Range("D32").Style = "Normal"

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0806'
$ws.Range("E33").Value = '  -3.47%  '
$ws.Range("D33").Style = "Normal"

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("B34").Value = 'LidoDAOToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D34").Value = '3.25'
$ws.Range("E34").Value = '  +3.42%  '
$ws.Range("D34").Style = "Normal"

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("B35").Value = 'WEMIXToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D35").Value = '2.61'
$ws.Range("E35").Value = '  -2.96%  '
$ws.Range("D35").Style = "Normal"

# Row 36
$ws.Range("E36").Value = '  +0.11%  '

# Row 37
$ws.Range("E37").Value = '  +2.70%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.78'
$ws.Range("E38").Value = '  -5.08%  '
$ws.Range("D38").Style = "Normal"

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '15.26'
$ws.Range("E39").Value = '  -3.78%  '
$ws.Range("D39").Style = "Normal"

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.43'
$ws.Range("E40").Value = '  -5.17%  '
$ws.Range("D40").Style = "Normal"

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.87'
$ws.Range("E41").Value = '  -3.81%  '
$ws.Range("D41").Style = "Normal"

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0300'
$ws.Range("E42").Value = '  -3.46%  '
$ws.Range("D42").Style = "Normal"

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.00'
$ws.Range("E43").Value = '  +0.08%  '
$ws.Range("D43").Style = "Normal"

# Row 44
$ws.Range("D44").Value = '1.745.95'
$ws.Range("E44").Value = '  +0.93%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '84.87'
$ws.Range("E45").Value = '  +5.00%  '
$ws.Range("D45").Style = "Normal"

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.189'
$ws.Range("E46").Value = '  -2.58%  '
$ws.Range("D46").Style = "Normal"

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '100.44'
$ws.Range("E47").Value = '  -0.96%  '
$ws.Range("D47").Style = "Normal"

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '4.97'
$ws.Range("E48").Value = '  -2.73%  '
$ws.Range("D48").Style = "Normal"

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("B49").Value = 'ordi'
$ws.Range("C49").Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range("D49").Value = '70.19'
$ws.Range("E49").Value = '  -4.18%  '
$ws.Range("D49").Style = "Normal"

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("B50").Value = 'FraxShare'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D50").Value = '8.11'
$ws.Range("E50").Value = '  +0.00%  '
$ws.Range("D50").Style = "Normal"

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("B51").Value = 'MultiversX'
$ws.Range("C51").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D51").Value = '54.37'
$ws.Range("E51").Value = '  -4.38%  '
$ws.Range("D51").Style = "Normal"
